$d = $word.ActiveDocument

# ============================================================
# Part 1: Insert new 'Knarot' section content after 'BILAGA 1 - Fridlysta arter'
# ============================================================

# Create all new empty paragraphs first (avoids paragraph style / run-formatting
# bleed-through from one new paragraph to the next).
for ($k = 0; $k -lt 13; $k++) {
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertParagraphAfter()
}

$baseIdx = $d.Paragraphs.Count - 12

# --- Paragraph 0 (style=Heading1) ---
$p0 = $d.Paragraphs($baseIdx)
$p0.Style = 'Heading1'
$p0.Range.Text = 'Knärot – ekologi samt krav på livsmiljön'

# --- Paragraph 1 (style=Normal) ---
$p1 = $d.Paragraphs($baseIdx + 1)
$p1.Style = 'Normal'
$p1.Range.Text = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'

# --- Paragraph 2 (style=Normal) ---
$p2 = $d.Paragraphs($baseIdx + 2)
$p2.Style = 'Normal'
$p2.Range.Text = 'Samuel Johnsons doktorsavhandling '
$ins = $p2.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“')
$ins.Font.Italic = $true
$ins = $p2.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ')
$ins = $p2.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ')
$ins.Font.Italic = $true
$ins = $p2.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Vidare ')
$ins = $p2.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$ins.Font.Italic = $true

# --- Paragraph 3 (style=Normal) ---
$p3 = $d.Paragraphs($baseIdx + 3)
$p3.Style = 'Normal'
$p3.Range.Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: '
$ins = $p3.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$ins.Font.Italic = $true

# --- Paragraph 4 (style=Normal) ---
$p4 = $d.Paragraphs($baseIdx + 4)
$p4.Style = 'Normal'
$p4.Range.Text = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'

# --- Paragraph 5 (style=Normal) ---
$p5 = $d.Paragraphs($baseIdx + 5)
$p5.Style = 'Normal'
$p5.Range.Text = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'

# --- Paragraph 6 (style=Heading2) ---
$p6 = $d.Paragraphs($baseIdx + 6)
$p6.Style = 'Heading2'
$p6.Range.Text = 'Referenser - knärot'

# --- Paragraph 7 (style=Normal) ---
$p7 = $d.Paragraphs($baseIdx + 7)
$p7.Style = 'Normal'
$p7.Range.Text = 'de Graaf M & Roberts M.R., 2009. '
$ins = $p7.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Short-term response of the herbaceous layer within leave patches after harvest. ')
$ins.Font.Italic = $true
$ins = $p7.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Forest Ecology and Management 257, 1014-1025')

# --- Paragraph 8 (style=Normal) ---
$p8 = $d.Paragraphs($baseIdx + 8)
$p8.Style = 'Normal'
$p8.Range.Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. '
$ins = $p8.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ')
$ins.Font.Italic = $true
$ins = $p8.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Ecological Applications, 22, 2049-2064 ')

# --- Paragraph 9 (style=Normal) ---
$p9 = $d.Paragraphs($baseIdx + 9)
$p9.Style = 'Normal'
$p9.Range.Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. '
$ins = $p9.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Interactive effects of drought and edge exposure on old-growth forest understory species. ')
$ins.Font.Italic = $true
$ins = $p9.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Landscape Ecology, 37, sid 1839-1853')

# --- Paragraph 10 (style=Normal) ---
$p10 = $d.Paragraphs($baseIdx + 10)
$p10.Style = 'Normal'
$p10.Range.Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. '
$ins = $p10.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Biological legacies buffer local species extinction after logging. ')
$ins.Font.Italic = $true
$ins = $p10.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Journal of Applied Ecology. 51, 53-62.')

# --- Paragraph 11 (style=Normal) ---
$p11 = $d.Paragraphs($baseIdx + 11)
$p11.Style = 'Normal'
$p11.Range.Text = 'Skogsstyrelsen, 2022. '
$ins = $p11.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Vägledning för hänsyn till knärot. ')
$ins.Font.Italic = $true
$ins = $p11.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')

# --- Paragraph 12 (style=Normal) ---
$p12 = $d.Paragraphs($baseIdx + 12)
$p12.Style = 'Normal'
$p12.Range.Text = 'SLU Artdatabanken, 2021. '
$ins = $p12.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('Artfaktablad. Naturvård – artfakta. ')
$ins.Font.Italic = $true
$ins = $p12.Range
$ins.SetRange($ins.End - 1, $ins.End - 1)
$ins.InsertAfter('SLU Artdatabanken, Uppsala ')

# ============================================================
# Part 2: Update the date in the first-page header from 2023-09-13 to 2023-09-15
# ============================================================
$sec = $d.Sections(1)
$firstPageHeader = $sec.Headers.Item(2)  # wdHeaderFooterFirstPage
$dateResult = $firstPageHeader.Range.Find.Execute("2023-09-13", $false, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)
Write-Host "Date replace result:" $dateResult
